# Update the "想去人数" (interested-count) figures in column F across the
# workbook's sheets, per the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2799
$ws.Range("F3").Value = 1127
$ws.Range("F4").Value = 20410
$ws.Range("F5").Value = 90
$ws.Range("F6").Value = 2519
$ws.Range("F7").Value = 779
$ws.Range("F8").Value = 613
$ws.Range("F9").Value = 476
$ws.Range("F10").Value = 729
$ws.Range("F11").Value = 266
$ws.Range("F15").Value = 96
$ws.Range("F16").Value = 492
$ws.Range("F17").Value = 176
$ws.Range("F18").Value = 237
$ws.Range("F21").Value = 111

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 24
$ws.Range("F5").Value = 309
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 119

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6071
$ws.Range("F3").Value = 676
$ws.Range("F4").Value = 646
$ws.Range("F5").Value = 1363
$ws.Range("F6").Value = 33

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6071
$ws.Range("F3").Value = 676
$ws.Range("F4").Value = 646
$ws.Range("F5").Value = 1363
$ws.Range("F6").Value = 2799
$ws.Range("F7").Value = 1127
$ws.Range("F8").Value = 20410
$ws.Range("F10").Value = 24
$ws.Range("F11").Value = 90
$ws.Range("F13").Value = 309
$ws.Range("F14").Value = 2519
$ws.Range("F15").Value = 779
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 613
$ws.Range("F19").Value = 477
$ws.Range("F20").Value = 729
$ws.Range("F21").Value = 266
$ws.Range("F27").Value = 391
$ws.Range("F28").Value = 96
$ws.Range("F31").Value = 492
$ws.Range("F33").Value = 176
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 237
$ws.Range("F36").Value = 119
$ws.Range("F37").Value = 119
$ws.Range("F48").Value = 111
